# Apply the edit described by the commit:
#   Storage: simple cyclic SoC - Formulation adapted from PyPSA
#   - Lack of high VRE nullifies storage
#   - Constant import price nullifies storage
#
# Concretely (as derived from the OOXML diff): a new data row is inserted
# at row 7 of Sheet1 ("enable_year" / "configuration" / 1990), pushing all
# the following rows down by one. The autofilter range, the hidden
# _FilterDatabase defined name and the sheet dimension all grow by one row
# accordingly, and the active selection ends up on H7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 7 (shifts rows 7.. down by one).
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new "enable_year" entry.
$ws.Range("A7").Value = "CHE"
$ws.Range("B7").Value = "conv_chp_coal"
$ws.Range("C7").Value = "enable_year"
$ws.Range("D7").Value = "configuration"
$ws.Range("G7").Value = 1990

# Grow the AutoFilter range from A5:L852 to A5:L853 to cover the new row.
$ws.AutoFilterMode = $false
$ws.Range("A5:L853").AutoFilter()

# Update the hidden _FilterDatabase defined name to match the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$5:`$L`$853"
    }
}

# Match the final active selection recorded in the saved workbook.
$ws.Range("H7").Select()
